$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 685 ("「ガザ地区」..." entry) entirely; Excel shifts rows 686:800 up by one.
$ws.Rows.Item(685).Delete()
